# Issue#30  Req 5.5  fix 7.3 - run happy path and make sure file is created
#
# Rewrites the three employee data rows on the "new sheet" worksheet from the
# old (Name/Surname/Division/...) sample data to the new Mentee/Mentor/
# BestMentor "happy path" sample data, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (was Tomasz Jurek / Information Systems Analyst / Corporate Services)
$ws.Range("B2").Value = "Mentee"
$ws.Range("C2").Value = "Developer"

# --- Row 3 (was Dariusz Łęcki / .Net Technical Lead / Business Consulting)
$ws.Range("B3").Value = "Mentor"
$ws.Range("C3").Value = "Developer"

# --- Row 4 (was Szymon Kaczmarczyk / Corporate PM & BA / Project Governance)
$ws.Range("B4").Value = "BestMentor"
$ws.Range("C4").Value = "Developer"

# --- Job family column (I) for all three rows
$ws.Range("I2").Value = "Project Development"
$ws.Range("I3").Value = "Project Development"
$ws.Range("I4").Value = "Project Development"

# --- Technology (G) / Position (H) columns for all three rows
$ws.Range("G2").Value = "Java"
$ws.Range("H2").Value = "Developer"
$ws.Range("G3").Value = "JavaScript"
$ws.Range("H3").Value = "Developer"
$ws.Range("G4").Value = "Java"
$ws.Range("H4").Value = "Developer"

# --- Division (E) / Grade (F) columns
$ws.Range("E2").Value = "Delivery"
$ws.Range("F2").Value = "L3"
$ws.Range("E3").Value = "Delivery"
$ws.Range("E4").Value = "Delivery"

# --- Office location (L) for row 4 moves from Warszawa to Łódź
$ws.Range("L4").Value = "Łódź"

# --- Move the sheet's active selection to H16
$ws.Range("H16").Select()
